$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C ("Förändrad") for rows 2 through 135 from 45186 to 45188
for ($r = 2; $r -le 135; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45186) {
        $cell.Value = 45188
    }
}
